$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("영화목록")

$ws.Range("D2").Value = "https://movie-phinf.pstatic.net/20220516_144/1652665409592Chvey_JPEG/movie_image.jpg"
$ws.Range("D3").Value = "https://movie-phinf.pstatic.net/20191024_143/1571893663418kwLN7_JPEG/movie_image.jpg"
$ws.Range("D4").Value = "https://movie-phinf.pstatic.net/20111224_165/13246577572754h14b_JPEG/movie_image.jpg"
$ws.Range("D5").Value = "https://movie-phinf.pstatic.net/20130204_279/1359954210596SuaVm_JPEG/movie_image.jpg"
$ws.Range("D6").Value = "https://movie-phinf.pstatic.net/20111223_87/13245909373833oAeh_JPEG/movie_image.jpg"
